# Scheduled-runner price/profit refresh across all Leve-profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose upstream market data changed.

$wb = $excel.ActiveWorkbook

function Set-Cell {
    param($ws, [string]$addr, $value)
    $ws.Range($addr).Value = $value
}

function Clear-Cell {
    param($ws, [string]$addr)
    $ws.Range($addr).ClearContents()
}

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-Cell $ws "H19" 1512.3334
Set-Cell $ws "I19" 1944.6666
Set-Cell $ws "J19" 1080
Set-Cell $ws "K19" 1944.6666
Set-Cell $ws "L19" 1080
Set-Cell $ws "M19" -1769.6666
Set-Cell $ws "N19" -1430

Set-Cell $ws "H33" 148.14285
Set-Cell $ws "I33" 148.14285
Set-Cell $ws "K33" 148.14285
Set-Cell $ws "M33" 80.85714999999999

Set-Cell $ws "H64" 3999.8333
Set-Cell $ws "I64" 0
Set-Cell $ws "J64" 3999.8333
Set-Cell $ws "K64" 0
Set-Cell $ws "L64" 3999.8333
Clear-Cell $ws "M64"
Set-Cell $ws "N64" -4495.8333

Set-Cell $ws "H67" 3999.8333
Set-Cell $ws "I67" 0
Set-Cell $ws "J67" 3999.8333
Set-Cell $ws "K67" 0
Set-Cell $ws "L67" 3999.8333
Clear-Cell $ws "M67"
Set-Cell $ws "N67" -5715.8333

Set-Cell $ws "H76" 5542
Set-Cell $ws "I76" 5948.5
Set-Cell $ws "K76" 5948.5
Set-Cell $ws "M76" -5633.5

Set-Cell $ws "H79" 5542
Set-Cell $ws "I79" 5948.5
Set-Cell $ws "K79" 5948.5
Set-Cell $ws "M79" -4856.5

Set-Cell $ws "H107" 0
Set-Cell $ws "I107" 0
Set-Cell $ws "K107" 0
Clear-Cell $ws "M107"

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-Cell $ws "H2" 1933.25
Set-Cell $ws "I2" 1933.25
Set-Cell $ws "K2" 1933.25
Set-Cell $ws "M2" -1820.25

Set-Cell $ws "H101" 0
Set-Cell $ws "J101" 0
Set-Cell $ws "L101" 0
Clear-Cell $ws "N101"

Set-Cell $ws "H116" 1933.25
Set-Cell $ws "I116" 1933.25
Set-Cell $ws "K116" 1933.25
Set-Cell $ws "M116" 360.75

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

Set-Cell $ws "H3" 1933.25
Set-Cell $ws "I3" 1933.25
Set-Cell $ws "K3" 1933.25
Set-Cell $ws "M3" -1819.25

Set-Cell $ws "H92" 49996.5
Set-Cell $ws "J92" 49996.5
Set-Cell $ws "L92" 49996.5
Set-Cell $ws "N92" -54988.5

Set-Cell $ws "H100" 20000
Set-Cell $ws "J100" 20000
Set-Cell $ws "L100" 20000
Set-Cell $ws "N100" -22164

Set-Cell $ws "H107" 4430.1816
Set-Cell $ws "I107" 3955.3333
Set-Cell $ws "K107" 3955.3333
Set-Cell $ws "M107" -2035.3333

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

Set-Cell $ws "H41" 20000
Set-Cell $ws "I41" 0
Set-Cell $ws "J41" 20000
Set-Cell $ws "K41" 0
Set-Cell $ws "L41" 20000
Clear-Cell $ws "M41"
Set-Cell $ws "N41" -20856

Set-Cell $ws "H43" 12660.2
Set-Cell $ws "J43" 12660.2
Set-Cell $ws "L43" 12660.2
Set-Cell $ws "N43" -13028.2

Set-Cell $ws "H58" 1604.7693
Set-Cell $ws "J58" 294
Set-Cell $ws "L58" 294
Set-Cell $ws "N58" -700

Set-Cell $ws "H88" 15814
Set-Cell $ws "J88" 15814
Set-Cell $ws "L88" 15814
Set-Cell $ws "N88" -16626

Set-Cell $ws "H91" 15814
Set-Cell $ws "J91" 15814
Set-Cell $ws "L91" 15814
Set-Cell $ws "N91" -18622

Set-Cell $ws "H101" 12660.2
Set-Cell $ws "J101" 12660.2
Set-Cell $ws "L101" 12660.2
Set-Cell $ws "N101" -19150.2

Set-Cell $ws "H107" 921.8125
Set-Cell $ws "I107" 939.2857
Set-Cell $ws "J107" 799.5
Set-Cell $ws "K107" 939.2857
Set-Cell $ws "L107" 799.5
Set-Cell $ws "M107" 980.7143
Set-Cell $ws "N107" -4639.5

Set-Cell $ws "H122" 1448
Set-Cell $ws "I122" 1448
Set-Cell $ws "K122" 4344
Set-Cell $ws "M122" -1894

Set-Cell $ws "H136" 1604.7693
Set-Cell $ws "J136" 294
Set-Cell $ws "L136" 882
Set-Cell $ws "N136" -5982

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Cell $ws "H12" 341.44446
Set-Cell $ws "J12" 373.8
Set-Cell $ws "L12" 1121.4
Set-Cell $ws "N12" -1467.4

Set-Cell $ws "H113" 499
Set-Cell $ws "I113" 499.4
Set-Cell $ws "J113" 497
Set-Cell $ws "K113" 1498.2
Set-Cell $ws "L113" 1491
Set-Cell $ws "M113" 671.8000000000002
Set-Cell $ws "N113" -5831

Set-Cell $ws "H117" 443.1111
Set-Cell $ws "J117" 515
Set-Cell $ws "L117" 1545
Set-Cell $ws "N117" -8429

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

Set-Cell $ws "H107" 1210.091
Set-Cell $ws "I107" 1210.091
Set-Cell $ws "K107" 1210.091
Set-Cell $ws "M107" 709.9090000000001

Set-Cell $ws "H123" 74498.2
Set-Cell $ws "J123" 74498.2
Set-Cell $ws "L123" 74498.2
Set-Cell $ws "N123" -79398.2

Set-Cell $ws "H126" 9998.200000000001
Set-Cell $ws "I126" 9998.333000000001
Set-Cell $ws "J126" 9998
Set-Cell $ws "K126" 29994.999
Set-Cell $ws "L126" 29994
Set-Cell $ws "M126" -27524.999
Set-Cell $ws "N126" -34934

Set-Cell $ws "H132" 5417
Set-Cell $ws "I132" 5417
Set-Cell $ws "K132" 16251
Set-Cell $ws "M132" -13721

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-Cell $ws "H16" 1712.375
Set-Cell $ws "I16" 1599.8572
Set-Cell $ws "K16" 1599.8572
Set-Cell $ws "M16" -1429.8572

Set-Cell $ws "H22" 4199.75
Set-Cell $ws "I22" 1000
Set-Cell $ws "J22" 5266.3335
Set-Cell $ws "K22" 1000
Set-Cell $ws "L22" 5266.3335
Set-Cell $ws "M22" -705
Set-Cell $ws "N22" -5856.3335

Set-Cell $ws "H27" 4199.75
Set-Cell $ws "I27" 1000
Set-Cell $ws "J27" 5266.3335
Set-Cell $ws "K27" 1000
Set-Cell $ws "L27" 5266.3335
Set-Cell $ws "M27" -893
Set-Cell $ws "N27" -5480.3335

Set-Cell $ws "H132" 14004
Set-Cell $ws "I132" 14004
Set-Cell $ws "J132" 0
Set-Cell $ws "K132" 42012
Set-Cell $ws "L132" 0
Set-Cell $ws "M132" -39482
Clear-Cell $ws "N132"

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

Set-Cell $ws "H95" 0
Set-Cell $ws "J95" 0
Set-Cell $ws "L95" 0
Clear-Cell $ws "N95"

Set-Cell $ws "H105" 25889.2
Set-Cell $ws "J105" 25889.2
Set-Cell $ws "L105" 25889.2
Set-Cell $ws "N105" -32877.2

Set-Cell $ws "H132" 3248.5
Set-Cell $ws "I132" 3248.5
Set-Cell $ws "K132" 9745.5
Set-Cell $ws "M132" -7215.5
